$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73162704"
$ws.Range("D16").Value = "WILSON FRIAS ALCALA"
$ws.Range("E16").Value = "2112"
$ws.Range("F16").Value = 36341
$ws.Range("G16").Value = 908526

# Row 17
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73574969"
$ws.Range("D17").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E17").Value = "2201"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 877803

# Row 18
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73574969"
$ws.Range("D18").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E18").Value = "2202"
$ws.Range("F18").Value = 36341
$ws.Range("G18").Value = 877803

# Row 19
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73574969"
$ws.Range("D19").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E19").Value = "2203"
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 877803

# Row 20
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1047456993"
$ws.Range("D20").Value = "ESTEBAN DE JESUS AHUMEDO BURGOS"
$ws.Range("E20").Value = "2204"
$ws.Range("F20").Value = 8000
$ws.Range("G20").Value = 1000000

# Row 21
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "73574969"
$ws.Range("D21").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E21").Value = "2204"
$ws.Range("F21").Value = 36341
$ws.Range("G21").Value = 877803

# Row 22
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "73574969"
$ws.Range("D22").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E22").Value = "2205"
$ws.Range("F22").Value = 36341
$ws.Range("G22").Value = 877803

# Row 23
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "73574969"
$ws.Range("D23").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E23").Value = "2206"
$ws.Range("F23").Value = 36341
$ws.Range("G23").Value = 877803

# Row 24
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "73574969"
$ws.Range("D24").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E24").Value = "2207"
$ws.Range("F24").Value = 36341
$ws.Range("G24").Value = 877803

# Row 25
$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "73574969"
$ws.Range("D25").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E25").Value = "2208"
$ws.Range("F25").Value = 36341
$ws.Range("G25").Value = 877803

# Row 26
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "73574969"
$ws.Range("D26").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E26").Value = "2209"
$ws.Range("F26").Value = 36341
$ws.Range("G26").Value = 877803

# Row 27
$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "73574969"
$ws.Range("D27").Value = "JOHN ALEXANDER PEREZ DE CASTRO"
$ws.Range("E27").Value = "2302"
$ws.Range("F27").Value = 46400
$ws.Range("G27").Value = 877803

# Row 28
$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1047402961"
$ws.Range("D28").Value = "JONATHAN MEZA BUSTAMANTE"
$ws.Range("E28").Value = "2304"
$ws.Range("F28").Value = 46400
$ws.Range("G28").Value = 1423500

# Row 29
$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "73119747"
$ws.Range("D29").Value = "PEDRO CLAVER CARABALLO OROZCO"
$ws.Range("E29").Value = "2304"
$ws.Range("F29").Value = 46400
$ws.Range("G29").Value = 1300000

# Row 30
$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "9296373"
$ws.Range("D30").Value = "RALLPH ANTONIO LLOREDA DORIA"
$ws.Range("E30").Value = "2304"
$ws.Range("F30").Value = 46400
$ws.Range("G30").Value = 1423500

# Row 31
$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "73167930"
$ws.Range("D31").Value = "OMAR DE JESUS TORRES CASTILLO"
$ws.Range("E31").Value = "2304"
$ws.Range("F31").Value = 46400
$ws.Range("G31").Value = 1423500
